$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.216.55"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.849.10"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7033"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07711"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3067"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.64"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07811"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "93.31"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.141"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "1.845.85"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6874"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.591"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008321"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "29.200.15"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.97"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").Value = "2.094.42"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.515"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9995"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1507"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.28"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.854"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.536"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.229"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.180"
$ws.Range("D31").ClearFormats()
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05122"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7878"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.894"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "1.322.33"
$ws.Range("E38").Value = "  +7.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01869"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.713"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9621"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.055"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.92"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.718"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "1.991.80"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5182"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.62"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.765"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.985"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.82%  "
